$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price-report row was recorded for "Superior Seedless" grapes
# (Región de O'Higgins market report, EE.UU. origin) and inserted as row 463,
# pushing every subsequent row (old 463..539) down by one (new 464..540).
$ws.Rows("463").Insert()

$ws.Range("A463").Value = 10
$ws.Range("B463").Value = "Vega Modelo de Temuco"
$ws.Range("C463").Value = "La Araucanía"
$ws.Range("D463").Value = 44505
$ws.Range("E463").Value = 9
$ws.Range("F463").Value = "Fruta"
$ws.Range("G463").Value = 100109
$ws.Range("H463").Value = "Uva"
$ws.Range("I463").Value = 100109001
$ws.Range("J463").Value = "Uva"
$ws.Range("K463").Value = "Superior Seedless"
$ws.Range("L463").Value = "Primera"
$ws.Range("M463").Value = 250
$ws.Range("N463").Value = 33000
$ws.Range("O463").Value = 34000
$ws.Range("P463").Value = 33400
$ws.Range("Q463").Value = "$/bandeja 8 kilos"
$ws.Range("R463").Value = "EE.UU."
$ws.Range("S463").Value = 4175
$ws.Range("T463").Value = 8
